$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(21, 35, 36)
foreach ($r in $rows) {
    $ws.Range("I$r").Value = "sv"
    $ws.Range("J$r").Value = "Statement-opinion"
}
